# Auto-generated edit script applying the cryptos.xlsx price/volume/name refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '86.760.05'
$ws.Range('E2').Value = '  +2.75%  '
$ws.Range('D3').Value = '3.244.60'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '210.05'
$ws.Range('E5').Value = '  -4.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '621.52'
$ws.Range('E6').Value = '  -2.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.360'
$ws.Range('E7').Value = '  +12.75%  '
$ws.Range('E8').Value = '  +13.29%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').Value = '3.240.30'
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.567'
$ws.Range('E11').Value = '  -5.35%  '
$ws.Range('E12').Value = '  +7.32%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000252'
$ws.Range('E13').Value = '  -9.20%  '
$ws.Range('D14').Value = '3.837.82'
$ws.Range('E14').Value = '  -1.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '33.63'
$ws.Range('E15').Value = '  +0.48%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.28'
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('D17').Value = '86.487.83'
$ws.Range('E17').Value = '  +2.58%  '
$ws.Range('D18').Value = '3.234.92'
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.08'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.89'
$ws.Range('E20').Value = '  -3.97%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '428.94'
$ws.Range('E21').Value = '  -4.87%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.79'
$ws.Range('E22').Value = '  -4.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.30'
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('E24').Value = '  -3.19%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.42'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '5.05'
$ws.Range('E26').Value = '  -4.17%  '
$ws.Range('D27').Value = '3.413.77'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '75.07'
$ws.Range('E28').Value = '  -3.59%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0000126'
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.173'
$ws.Range('E31').Value = '  +10.50%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '8.69'
$ws.Range('E33').Value = '  -5.86%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '538.38'
$ws.Range('E34').Value = '  -5.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.39'
$ws.Range('E35').Value = '  -7.99%  '
$ws.Range('E36').Value = '  -4.91%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.92'
$ws.Range('E37').Value = '  +10.85%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.135'
$ws.Range('E38').Value = '  -12.14%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '22.17'
$ws.Range('E39').Value = '  -4.87%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '21.61'
$ws.Range('E41').Value = '  +3.17%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.98'
$ws.Range('E42').Value = '  -4.07%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.387'
$ws.Range('E43').Value = '  -6.35%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.89'
$ws.Range('E45').Value = '  -5.53%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '153.99'
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '177.33'
$ws.Range('E47').Value = '  -7.16%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '44.27'
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.27'
$ws.Range('E49').Value = '  -4.67%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.17'
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('E51').Value = '  +9.95%  '
